# Update "想去人数" (F column) and occasionally "最低票价" (G column)
# figures across the "展览" (sheet1), "演出" (sheet2), "本地生活" (sheet3)
# and "全部类型" (sheet4) worksheets, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# ---- 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 11191
$ws1.Cells.Item(5, 6).Value = 1256
$ws1.Cells.Item(6, 6).Value = 1130
$ws1.Cells.Item(7, 6).Value = 874
$ws1.Cells.Item(11, 6).Value = 160
$ws1.Cells.Item(12, 6).Value = 930
$ws1.Cells.Item(13, 6).Value = 2164
$ws1.Cells.Item(15, 6).Value = 1067
$ws1.Cells.Item(19, 6).Value = 969
$ws1.Cells.Item(21, 6).Value = 274
$ws1.Cells.Item(23, 6).Value = 667
$ws1.Cells.Item(24, 6).Value = 694
$ws1.Cells.Item(25, 6).Value = 142
$ws1.Cells.Item(26, 6).Value = 381
$ws1.Cells.Item(27, 6).Value = 1035
$ws1.Cells.Item(29, 6).Value = 416
$ws1.Cells.Item(29, 7).Value = 40.5
$ws1.Cells.Item(33, 6).Value = 258
$ws1.Cells.Item(34, 6).Value = 608
$ws1.Cells.Item(35, 6).Value = 2228
$ws1.Cells.Item(35, 7).Value = 67.5
$ws1.Cells.Item(36, 6).Value = 415
$ws1.Cells.Item(38, 6).Value = 1484
$ws1.Cells.Item(43, 6).Value = 50
$ws1.Cells.Item(45, 6).Value = 90

# ---- 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 208

# ---- 本地生活 ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 2204
$ws3.Cells.Item(4, 6).Value = 610

# ---- 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 2204
$ws4.Cells.Item(5, 6).Value = 1256
$ws4.Cells.Item(6, 6).Value = 610
$ws4.Cells.Item(7, 6).Value = 1130
$ws4.Cells.Item(8, 6).Value = 874
$ws4.Cells.Item(9, 6).Value = 208
$ws4.Cells.Item(14, 6).Value = 930
$ws4.Cells.Item(15, 6).Value = 2164
$ws4.Cells.Item(17, 6).Value = 1067
$ws4.Cells.Item(21, 6).Value = 969
$ws4.Cells.Item(22, 6).Value = 274
$ws4.Cells.Item(25, 6).Value = 667
$ws4.Cells.Item(26, 6).Value = 694
$ws4.Cells.Item(27, 6).Value = 142
$ws4.Cells.Item(28, 6).Value = 381
$ws4.Cells.Item(29, 6).Value = 1035
$ws4.Cells.Item(31, 6).Value = 416
$ws4.Cells.Item(31, 7).Value = 40.5
$ws4.Cells.Item(35, 6).Value = 258
$ws4.Cells.Item(36, 6).Value = 2229
$ws4.Cells.Item(36, 7).Value = 67.5
$ws4.Cells.Item(38, 6).Value = 415
$ws4.Cells.Item(40, 6).Value = 1484
$ws4.Cells.Item(44, 6).Value = 50
$ws4.Cells.Item(45, 6).Value = 90
